# 강민재 enum LaneBuffType 추가
# Adds a new "LaneBuffType" enum table (rows 123-151) and logs the change
# in the revision-history block (row 33, columns I:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New revision-log row, right below the previous last entry (row 32).
#    Columns: I=date, J=author, K=modified enum, L=action.
# ---------------------------------------------------------------------
$ws.Range("I32:L32").Copy()
$ws.Range("I33:L33").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I33").Value = "09.06 13:49"
$ws.Range("J33").Value = "강민재"
$ws.Range("K33").Value = "LaneBuffType"
$ws.Range("L33").Value = "추가"

# ---------------------------------------------------------------------
# 2) New enum table: LaneBuffType, starting at row 123
#    (two blank rows after the previous ElementalType table which ends
#    at row 120).
# ---------------------------------------------------------------------

# Header row - copy formatting from the matching header of the
# ElementalType table (row 115: EnumName / Typename / TypeValue).
$ws.Range("A115:C115").Copy()
$ws.Range("A123:C123").PasteSpecial(-4122)

$ws.Range("A123").Value = "EnumName"
$ws.Range("B123").Value = "Typename"
$ws.Range("C123").Value = "TypeValue"

# Data rows - copy formatting from an existing enum block that uses the
# same style pattern (A: s=1, B: s=39, C: s=1), e.g. rows 77-89.
$ws.Range("A77:C89").Copy()
$ws.Range("A124:C136").PasteSpecial(-4122)
$ws.Range("A77:C89").Copy()
$ws.Range("A137:C149").PasteSpecial(-4122)
$ws.Range("A77:C78").Copy()
$ws.Range("A150:C151").PasteSpecial(-4122)

$ws.Range("A124").Value = "LaneBuffType"

$members = @(
    "AllHealHp",
    "AllDealHp",
    "AllHealShield",
    "AllDealShield",
    "AllBuffArmor",
    "AllNurfArmor",
    "WaterHealHp",
    "GroundHealHp",
    "FireHealHp",
    "ElectricHealHp",
    "WaterDealHp",
    "GroundDealHp",
    "FireHealHp",
    "ElectricDealHp",
    "WaterBuffArmor",
    "GroundBuffArmor",
    "FireBuffArmor",
    "ElectricBuffArmor",
    "WaterNurfArmor",
    "GroundNurfArmor",
    "FireNurfArmor",
    "ElectricNurfArmor",
    "ArrowBuff",
    "SlowBuff",
    "PoisonBuff",
    "FlameBuff",
    "LazerBuff",
    "MissileBuff"
)

for ($i = 0; $i -lt $members.Length; $i++) {
    $r = 124 + $i
    $ws.Cells.Item($r, 2).Value = $members[$i]
    $ws.Cells.Item($r, 3).Value = $i + 1
}

# ---------------------------------------------------------------------
# 3) Selection / viewport, matching where the author left off editing.
# ---------------------------------------------------------------------
$ws.Range("L33").Select()
$excel.ActiveWindow.ScrollRow = 19
